# Update market/profit data columns (currentAveragePrice, currentAveragePriceNQ/HQ,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) across all leve sheets, as produced by the
# scheduled market-data refresh runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6649.7
$ws.Range("I18").Value = 7277.4443
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 7277.4443
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = -6993.4443
$ws.Range("N18").Value = -1568
$ws.Range("H53").Value = 635.875
$ws.Range("I53").Value = 532
$ws.Range("J53").Value = 698.2
$ws.Range("K53").Value = 532
$ws.Range("L53").Value = 698.2
$ws.Range("M53").Value = 105
$ws.Range("N53").Value = -1972.2
$ws.Range("H62").Value = 4158.8
$ws.Range("I62").Value = 4158.8
$ws.Range("K62").Value = 4158.8
$ws.Range("M62").Value = -3534.8
$ws.Range("H65").Value = 4158.8
$ws.Range("I65").Value = 4158.8
$ws.Range("K65").Value = 20794
$ws.Range("M65").Value = -17674
$ws.Range("H92").Value = 866.86664
$ws.Range("I92").Value = 901.2857
$ws.Range("K92").Value = 901.2857
$ws.Range("M92").Value = 346.7143
$ws.Range("H98").Value = 22360
$ws.Range("I98").Value = 29666.105
$ws.Range("K98").Value = 29666.105
$ws.Range("M98").Value = -28168.105
$ws.Range("H122").Value = 22360
$ws.Range("I122").Value = 29666.105
$ws.Range("K122").Value = 88998.315
$ws.Range("M122").Value = -86548.315
$ws.Range("H132").Value = 3828.8
$ws.Range("I132").Value = 3727.257
$ws.Range("K132").Value = 11181.771
$ws.Range("M132").Value = -8651.771000000001
$ws.Range("H135").Value = 4460.1875
$ws.Range("I135").Value = 6014.2856
$ws.Range("J135").Value = 1493.2727
$ws.Range("K135").Value = 54128.5704
$ws.Range("L135").Value = 13439.4543
$ws.Range("M135").Value = -51593.5704
$ws.Range("N135").Value = -18509.4543

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 58250.406
$ws.Range("I45").Value = 129324.69
$ws.Range("J45").Value = 4098.5713
$ws.Range("K45").Value = 129324.69
$ws.Range("L45").Value = 4098.5713
$ws.Range("M45").Value = -128947.69
$ws.Range("N45").Value = -4852.5713
$ws.Range("H61").Value = 7584.1714
$ws.Range("I61").Value = 10017.667
$ws.Range("K61").Value = 10017.667
$ws.Range("M61").Value = -9805.666999999999
$ws.Range("H97").Value = 7696691
$ws.Range("I97").Value = 6716.6875
$ws.Range("J97").Value = 20000650
$ws.Range("K97").Value = 6716.6875
$ws.Range("L97").Value = 20000650
$ws.Range("M97").Value = -6220.6875
$ws.Range("N97").Value = -20001642
$ws.Range("H102").Value = 7744.3823
$ws.Range("I102").Value = 9802.32
$ws.Range("K102").Value = 9802.32
$ws.Range("M102").Value = -8180.32
$ws.Range("H132").Value = 2892.4614
$ws.Range("I132").Value = 2838.1
$ws.Range("J132").Value = 3073.6667
$ws.Range("K132").Value = 8514.299999999999
$ws.Range("L132").Value = 9221.000100000001
$ws.Range("M132").Value = -5984.299999999999
$ws.Range("N132").Value = -14281.0001
$ws.Range("H136").Value = 7584.1714
$ws.Range("I136").Value = 10017.667
$ws.Range("K136").Value = 30053.001
$ws.Range("M136").Value = -27503.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 8968.806
$ws.Range("I94").Value = 11661.32
$ws.Range("K94").Value = 11661.32
$ws.Range("M94").Value = -11210.32
$ws.Range("H99").Value = 9407.706
$ws.Range("I99").Value = 9780.75
$ws.Range("J99").Value = 7666.8335
$ws.Range("K99").Value = 9780.75
$ws.Range("L99").Value = 7666.8335
$ws.Range("M99").Value = -8282.75
$ws.Range("N99").Value = -10662.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1953.7
$ws.Range("I22").Value = 1923.3334
$ws.Range("K22").Value = 1923.3334
$ws.Range("M22").Value = -1573.3334
$ws.Range("H31").Value = 13225.462
$ws.Range("I31").Value = 41998
$ws.Range("J31").Value = 4593.7
$ws.Range("K31").Value = 41998
$ws.Range("L31").Value = 4593.7
$ws.Range("M31").Value = -41703
$ws.Range("N31").Value = -5183.7
$ws.Range("H34").Value = 13225.462
$ws.Range("I34").Value = 41998
$ws.Range("J34").Value = 4593.7
$ws.Range("K34").Value = 41998
$ws.Range("L34").Value = 4593.7
$ws.Range("M34").Value = -41796
$ws.Range("N34").Value = -4997.7
$ws.Range("H86").Value = 13999.75
$ws.Range("I86").Value = 9499.5
$ws.Range("K86").Value = 9499.5
$ws.Range("M86").Value = -8376.5
$ws.Range("H89").Value = 13999.75
$ws.Range("I89").Value = 9499.5
$ws.Range("K89").Value = 47497.5
$ws.Range("M89").Value = -41881.5
$ws.Range("H107").Value = 4190.4546
$ws.Range("I107").Value = 5139.5386
$ws.Range("K107").Value = 5139.5386
$ws.Range("M107").Value = -3219.5386
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 4638.0527
$ws.Range("I134").Value = 5483.893
$ws.Range("J134").Value = 2269.7
$ws.Range("K134").Value = 16451.679
$ws.Range("L134").Value = 6809.099999999999
$ws.Range("M134").Value = -13916.679
$ws.Range("N134").Value = -11879.1
$ws.Range("H137").Value = 41333.332
$ws.Range("J137").Value = 49500
$ws.Range("L137").Value = 49500
$ws.Range("N137").Value = -59700
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 18.416666
$ws.Range("J12").Value = 7.5555553
$ws.Range("L12").Value = 22.6666659
$ws.Range("N12").Value = -368.6666659
$ws.Range("H15").Value = 403.75
$ws.Range("I15").Value = 115
$ws.Range("K15").Value = 345
$ws.Range("M15").Value = -205
$ws.Range("H133").Value = 14742.375
$ws.Range("I133").Value = 5983
$ws.Range("K133").Value = 17949
$ws.Range("M133").Value = -12889
$ws.Range("H138").Value = 925.2857
$ws.Range("I138").Value = 940.6667
$ws.Range("K138").Value = 2822.0001
$ws.Range("M138").Value = 2317.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 19832.666
$ws.Range("J80").Value = 4750
$ws.Range("L80").Value = 4750
$ws.Range("N80").Value = -6746
$ws.Range("H83").Value = 19832.666
$ws.Range("J83").Value = 4750
$ws.Range("L83").Value = 23750
$ws.Range("N83").Value = -33734
$ws.Range("H132").Value = 4060.5745
$ws.Range("I132").Value = 4456.9443
$ws.Range("J132").Value = 2763.3635
$ws.Range("K132").Value = 13370.8329
$ws.Range("L132").Value = 8290.0905
$ws.Range("M132").Value = -10840.8329
$ws.Range("N132").Value = -13350.0905
$ws.Range("H137").Value = 45875
$ws.Range("J137").Value = 51750
$ws.Range("L137").Value = 51750
$ws.Range("N137").Value = -61950
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3725.6
$ws.Range("I16").Value = 3682.25
$ws.Range("J16").Value = 3899
$ws.Range("K16").Value = 3682.25
$ws.Range("L16").Value = 3899
$ws.Range("M16").Value = -3512.25
$ws.Range("N16").Value = -4239
$ws.Range("H46").Value = 3874245.8
$ws.Range("I46").Value = 861.25
$ws.Range("J46").Value = 6972953.5
$ws.Range("K46").Value = 861.25
$ws.Range("L46").Value = 6972953.5
$ws.Range("M46").Value = -673.25
$ws.Range("N46").Value = -6973329.5
$ws.Range("H55").Value = 2167.25
$ws.Range("I55").Value = 267.8
$ws.Range("J55").Value = 5333
$ws.Range("K55").Value = 267.8
$ws.Range("L55").Value = 5333
$ws.Range("M55").Value = -94.80000000000001
$ws.Range("N55").Value = -5679
$ws.Range("H68").Value = 5699
$ws.Range("I68").Value = 2118.4
$ws.Range("J68").Value = 11666.667
$ws.Range("K68").Value = 2118.4
$ws.Range("L68").Value = 11666.667
$ws.Range("M68").Value = -1369.4
$ws.Range("N68").Value = -13164.667
$ws.Range("H71").Value = 5699
$ws.Range("I71").Value = 2118.4
$ws.Range("J71").Value = 11666.667
$ws.Range("K71").Value = 10592
$ws.Range("L71").Value = 58333.335
$ws.Range("M71").Value = -6848
$ws.Range("N71").Value = -65821.33499999999
$ws.Range("H93").Value = 3595.8276
$ws.Range("I93").Value = 4191.7827
$ws.Range("J93").Value = 1311.3334
$ws.Range("K93").Value = 4191.7827
$ws.Range("L93").Value = 1311.3334
$ws.Range("M93").Value = -2943.7827
$ws.Range("N93").Value = -3807.3334
$ws.Range("H122").Value = 4801.054
$ws.Range("I122").Value = 4580.357
$ws.Range("K122").Value = 13741.071
$ws.Range("M122").Value = -11291.071
$ws.Range("H136").Value = 4173.4814
$ws.Range("I136").Value = 3206.6
$ws.Range("K136").Value = 9619.799999999999
$ws.Range("M136").Value = -7069.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 18127.955
$ws.Range("J100").Value = 53577.832
$ws.Range("L100").Value = 107155.664
$ws.Range("N100").Value = -108237.664
$ws.Range("H132").Value = 9748.257
$ws.Range("I132").Value = 14051.695
$ws.Range("J132").Value = 3562.0625
$ws.Range("K132").Value = 42155.085
$ws.Range("L132").Value = 10686.1875
$ws.Range("M132").Value = -39625.085
$ws.Range("N132").Value = -15746.1875
$ws.Range("H136").Value = 308401.4
$ws.Range("I136").Value = 354484.12
$ws.Range("J136").Value = 18738.572
$ws.Range("K136").Value = 1063452.36
$ws.Range("L136").Value = 56215.716
$ws.Range("M136").Value = -1060902.36
$ws.Range("N136").Value = -61315.716
$ws.Range("H140").Value = 101997.5
$ws.Range("J140").Value = 101997.5
$ws.Range("L140").Value = 101997.5
$ws.Range("N140").Value = -112357.5
